$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos Lab7")

# --- Fill in the measured PROBING results (rows 3-5: load factor 0.3/0.5/0.8) ---
$ws.Range("B3").Value = 1324030.27
$ws.Range("C3").Value = 48244.205999999998
$ws.Range("B4").Value = 1324034.3259999999
$ws.Range("C4").Value = 48030.326999999997
$ws.Range("B5").Value = 1324036.5530000001
$ws.Range("C5").Value = 47987.228000000003

# --- Fill in the measured CHAINING results (rows 10-12: load factor 2/4/6) ---
$ws.Range("B10").Value = 1324054.6359999999
$ws.Range("C10").Value = 49789.432000000001
$ws.Range("B11").Value = 1324055.2390000001
$ws.Range("C11").Value = 51158.557000000001
$ws.Range("B12").Value = 1324055.7949999999
$ws.Range("C12").Value = 54357.762000000002

# --- Now that there is real data, format it as thousands-separated numbers ---
$ws.Range("B3:C5").NumberFormat = "#,##0"
$ws.Range("B10:C12").NumberFormat = "#,##0"

# --- Header rows go back to the sheet's default height ---
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(9).AutoFit()

# --- Leave the sheet with the same view state as the final save ---
$ws.Range("D10").Select() | Out-Null
$excel.ActiveWindow.Zoom = 163
